$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C8").Value = 181345
$ws.Range("E8").Value = 649878687
$ws.Range("C10").Value = 278183
$ws.Range("E10").Value = 1751972720
$ws.Range("C17").Value = 134739
$ws.Range("E17").Value = 296786936
$ws.Range("C69").Value = 20735
$ws.Range("E69").Value = 62161544
$ws.Range("C81").Value = 26156
$ws.Range("D81").Value = 2827
$ws.Range("E81").Value = 165028835
$ws.Range("C99").Value = 136566
$ws.Range("E99").Value = 863015376
$ws.Range("C110").Value = 16863
$ws.Range("E110").Value = 25926995
$ws.Range("C111").Value = 6004
$ws.Range("E111").Value = 12042753
$ws.Range("C114").Value = 7488
$ws.Range("E114").Value = 11583350
$ws.Range("C115").Value = 17533
$ws.Range("E115").Value = 38563671
$ws.Range("C117").Value = 19686
$ws.Range("E117").Value = 56340101
$ws.Range("C120").Value = 2329
$ws.Range("E120").Value = 4371451
$ws.Range("C125").Value = 4590
$ws.Range("E125").Value = 13132883
$ws.Range("C126").Value = 5639
$ws.Range("E126").Value = 8168120
$ws.Range("C134").Value = 5662
$ws.Range("E134").Value = 17023521
$ws.Range("C150").Value = 95005
$ws.Range("E150").Value = 278392604
$ws.Range("C152").Value = 126037
$ws.Range("E152").Value = 715735940
$ws.Range("C162").Value = 62067
$ws.Range("E162").Value = 113487855
$ws.Range("C168").Value = 284884
$ws.Range("E168").Value = 1207126425
$ws.Range("C170").Value = 367189
$ws.Range("E170").Value = 2842709840
$ws.Range("C171").Value = 115074
$ws.Range("E171").Value = 444080453
$ws.Range("C174").Value = 357109
$ws.Range("E174").Value = 1015502936
$ws.Range("E175").Value = 809382495
$ws.Range("C177").Value = 96738
$ws.Range("E177").Value = 174252127
$ws.Range("C179").Value = 235613
$ws.Range("E179").Value = 811458870
$ws.Range("C180").Value = 141436
$ws.Range("E180").Value = 339794168
$ws.Range("C186").Value = 21924
$ws.Range("E186").Value = 39902432
$ws.Range("C188").Value = 19677
$ws.Range("E188").Value = 65875236
$ws.Range("C193").Value = 5341
$ws.Range("E193").Value = 27521411
$ws.Range("C196").Value = 7391
$ws.Range("E196").Value = 20454086
$ws.Range("C198").Value = 4507
$ws.Range("E198").Value = 5999037
$ws.Range("C203").Value = 13088
$ws.Range("E203").Value = 32905556
$ws.Range("C204").Value = 4747
$ws.Range("E204").Value = 11625430
$ws.Range("C205").Value = 11110
$ws.Range("E205").Value = 43922283
$ws.Range("C209").Value = 5361
$ws.Range("E209").Value = 12205320
$ws.Range("C213").Value = 3625
$ws.Range("E213").Value = 11062086
$ws.Range("C239").Value = 84899
$ws.Range("E239").Value = 500257654
$ws.Range("C257").Value = 182542
$ws.Range("E257").Value = 1063631676
$ws.Range("C258").Value = 15135
$ws.Range("E258").Value = 40650535
$ws.Range("C267").Value = 84971
$ws.Range("E267").Value = 156507125
$ws.Range("C283").Value = 60838
$ws.Range("E283").Value = 101611499
$ws.Range("C295").Value = 91330
$ws.Range("D295").Value = 9956
$ws.Range("E295").Value = 552903537
$ws.Range("C313").Value = 220629
$ws.Range("E313").Value = 1370587363
$ws.Range("C323").Value = 94722
$ws.Range("E323").Value = 178789585
